$d = $word.ActiveDocument

$pairs = @(
    @("686÷2=", "767÷5="),
    @("278÷7=", "998÷2="),
    @("660÷4=", "594÷9="),
    @("279÷4=", "347÷3="),
    @("399÷9=", "569÷9="),
    @("401÷9=", "766÷4="),
    @("496÷7=", "895÷8="),
    @("694÷7=", "888÷3="),
    @("430÷7=", "623÷3="),
    @("524÷8=", "776÷5="),
    @("169÷2=", "774÷3="),
    @("737÷9=", "841÷6="),
    @("186÷4=", "163÷6="),
    @("430÷2=", "370÷9="),
    @("314÷9=", "945÷4="),
    @("847÷7=", "367÷2="),
    @("808÷4=", "838÷5="),
    @("318÷6=", "773÷3="),
    @("639÷2=", "110÷9="),
    @("976÷7=", "197÷8="),
    @("568÷5=", "429÷6="),
    @("622÷2=", "337÷4="),
    @("833÷8=", "394÷7="),
    @("630÷4=", "208÷8="),
    @("692÷5=", "126÷5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
